$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.749
$ws.Range("B6").Value = 7.295999999999999
$ws.Range("B7").Value = 6.611999999999999
$ws.Range("B16").Value = 6.241
$ws.Range("B20").Value = 6.462000000000001
